$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerated sval data (filter save games) for rows 2-6, columns B:G.
# Note: the engine's expression parser does not accept scientific-notation
# numeric literals (e.g. 1.23e-05), so those are expressed as a division.

$ws.Cells.Item(2, 2).Value = 1.19090754144846 / 100000
$ws.Cells.Item(2, 3).Value = 0.002658071450198252
$ws.Cells.Item(2, 4).Value = 0.1496068669990043
$ws.Cells.Item(2, 5).Value = 13.86384647080068
$ws.Cells.Item(2, 6).Value = 0
$ws.Cells.Item(2, 7).Value = 14.0161233183253

$ws.Cells.Item(3, 2).Value = 0.6545652718822623
$ws.Cells.Item(3, 3).Value = 1.626987699542094
$ws.Cells.Item(3, 4).Value = 0.7210945179870265
$ws.Cells.Item(3, 5).Value = 0.5333859586016987
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 3.536033448013082

$ws.Cells.Item(4, 2).Value = 9.552326474482342 / 100000
$ws.Cells.Item(4, 3).Value = 6.194867796516235 / 10000000
$ws.Cells.Item(4, 4).Value = 0.1496068669990043
$ws.Cells.Item(4, 5).Value = 13.86384647080068
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value = 14.01354948055121

$ws.Cells.Item(5, 2).Value = 0.6545652718822623
$ws.Cells.Item(5, 3).Value = 0.3048912486333797
$ws.Cells.Item(5, 4).Value = 0.7210945179870265
$ws.Cells.Item(5, 5).Value = 0.5333859586016987
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 2.213936997104367

$ws.Cells.Item(6, 2).Value = 3.272327238179451
$ws.Cells.Item(6, 3).Value = 9.983522426115931
$ws.Cells.Item(6, 4).Value = 3.223369029078222
$ws.Cells.Item(6, 5).Value = 13.86384647080068
$ws.Cells.Item(6, 6).Value = 0
$ws.Cells.Item(6, 7).Value = 30.34306516417429
